# Regenerate merged AHB files
#
# 1. Rename the diff-header columns:
#      *_old -> *_FV2410   (columns A-J)
#      *_new -> *_FV2504   (columns L-U)
#    (column K, "diff", is unchanged)
# 2. Turn the A1:U65 range into a native Excel Table ("Table1") so the
#    renamed headers become the table's column names.
# 3. Freeze the header row (pane split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------
$headerNames = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerNames[$i]
}

# --- 2. Create the table over the used range ------------------------------
$tableRange = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the top (header) row ---------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
